$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column values (column G), replacing the old Strike# values
$kValues = @(
    1, 1, 1, 1, 2, 0, 2, 3, 1, 0, 2, 2, 1, 0, 2, 0, 2, 1, 0, 2, 4, 0, 1, 1, 2, 2, 0, 0, 2, 0, 3, 1, 1, 0, 0, 2, 1, 0, 2, 2, 2, 1, 1, 2, 2, 0, 3, 1, 1, 1, 2, 2, 1, 1, 0
)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
